$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.894.39"
$ws.Range("E2").Value = "  -2.29%  "
$ws.Range("D3").Value = "3.409.77"
$ws.Range("E3").Value = "  -2.95%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.30"
$ws.Range("E5").Value = "  -2.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "126.26"
$ws.Range("E6").Value = "  -6.21%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "3.410.07"
$ws.Range("E8").Value = "  -2.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.476"
$ws.Range("E9").Value = "  -2.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.40"
$ws.Range("E10").Value = "  -3.04%  "
$ws.Range("E11").Value = "  -2.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.378"
$ws.Range("E12").Value = "  -2.69%  "
$ws.Range("D13").Value = "3.987.36"
$ws.Range("E13").Value = "  -3.05%  "
$ws.Range("E14").Value = "  -0.58%  "
$ws.Range("D15").Value = "3.408.55"
$ws.Range("E15").Value = "  -3.00%  "
$ws.Range("E16").Value = "  -4.39%  "
$ws.Range("D17").Value = "62.944.67"
$ws.Range("E17").Value = "  -2.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "24.87"
$ws.Range("E18").Value = "  -3.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.54"
$ws.Range("E19").Value = "  -4.57%  "
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("E21").Value = "  -2.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "376.67"
$ws.Range("E22").Value = "  -4.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.560"
$ws.Range("E23").Value = "  -2.88%  "
$ws.Range("D24").Value = "3.546.34"
$ws.Range("E24").Value = "  -2.95%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -7.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("E29").Value = "  -5.45%  "
$ws.Range("E30").Value = "  -4.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.87"
$ws.Range("E31").Value = "  -5.15%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.40"
$ws.Range("E32").Value = "  -4.62%  "
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.151"
$ws.Range("E33").Value = "  -3.60%  "
$ws.Range("B34").Value = "RenzoRestakedETH"
$ws.Range("C34").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D34").Value = "3.440.68"
$ws.Range("E34").Value = "  -2.93%  "
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "22.80"
$ws.Range("E36").Value = "  -2.60%  "
$ws.Range("E37").Value = "  -0.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "165.74"
$ws.Range("E38").Value = "  -0.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.73"
$ws.Range("E39").Value = "  -3.36%  "
$ws.Range("E40").Value = "  -4.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0759"
$ws.Range("E41").Value = "  -3.89%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.777"
$ws.Range("E43").Value = "  -4.35%  "
$ws.Range("E44").Value = "  -1.64%  "
$ws.Range("E45").Value = "  -3.63%  "
$ws.Range("E46").Value = "  -5.22%  "
$ws.Range("E47").Value = "  -10.73%  "
$ws.Range("E48").Value = "  -8.02%  "
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("D50").Value = "2.259.80"
$ws.Range("E50").Value = "  -5.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.858"
$ws.Range("E51").Value = "  -4.52%  "
